$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 180; this pushes the existing rows
# 180..234 down to 181..235 (values, styles and formats shift with them).
$ws.Rows.Item(180).Insert()

# Populate the newly inserted row 180 with the new data point.
$ws.Cells.Item(180, 1).Value = 10
$ws.Cells.Item(180, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(180, 3).Value = "La Araucanía"
$ws.Cells.Item(180, 4).Value = 45215
$ws.Cells.Item(180, 5).Value = 9
$ws.Cells.Item(180, 6).Value = 100112031
$ws.Cells.Item(180, 7).Value = "Poroto verde"
$ws.Cells.Item(180, 8).Value = "Sin especificar"
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 65
$ws.Cells.Item(180, 11).Value = 1500
$ws.Cells.Item(180, 12).Value = 1500
$ws.Cells.Item(180, 13).Value = 1500
$ws.Cells.Item(180, 14).Value = "$/kilo"
$ws.Cells.Item(180, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(180, 16).Value = 1500
$ws.Cells.Item(180, 17).Value = 1
$ws.Cells.Item(180, 18).Value = "Hortaliza"
